{"js": "// Update the multiplication practice table with a new set of problems.\n// Each filled row (0-indexed: 0, 4, 9, 14, 19) gets all five of its\n// cell values replaced with new \"NNN\u00d7N=\" expressions. We replace the\n// text via each cell paragraph's own Range so the existing run/paragraph\n// formatting (font, size, justification) is preserved, and we address\n// cells positionally (row, column) so that values which happen to\n// collide with other rows' old/new text cannot cross-contaminate\n// each other.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"601\u00d75=\" },\n  { row: 0, col: 1, text: \"610\u00d73=\" },\n  { row: 0, col: 2, text: \"568\u00d77=\" },\n  { row: 0, col: 3, text: \"945\u00d78=\" },\n  { row: 0, col: 4, text: \"770\u00d72=\" },\n\n  { row: 4, col: 0, text: \"740\u00d78=\" },\n  { row: 4, col: 1, text: \"762\u00d78=\" },\n  { row: 4, col: 2, text: \"224\u00d78=\" },\n  { row: 4, col: 3, text: \"906\u00d73=\" },\n  { row: 4, col: 4, text: \"138\u00d76=\" },\n\n  { row: 9, col: 0, text: \"762\u00d76=\" },\n  { row: 9, col: 1, text: \"960\u00d77=\" },\n  { row: 9, col: 2, text: \"434\u00d78=\" },\n  { row: 9, col: 3, text: \"105\u00d74=\" },\n  { row: 9, col: 4, text: \"513\u00d78=\" },\n\n  { row: 14, col: 0, text: \"737\u00d78=\" },\n  { row: 14, col: 1, text: \"346\u00d76=\" },\n  { row: 14, col: 2, text: \"737\u00d73=\" },\n  { row: 14, col: 3, text: \"742\u00d78=\" },\n  { row: 14, col: 4, text: \"725\u00d77=\" },\n\n  { row: 19, col: 0, text: \"816\u00d78=\" },\n  { row: 19, col: 1, text: \"279\u00d78=\" },\n  { row: 19, col: 2, text: \"701\u00d72=\" },\n  { row: 19, col: 3, text: \"934\u00d76=\" },\n  { row: 19, col: 4, text: \"652\u00d79=\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const range = paragraphs.items[0].getRange();\n  range.insertText(u.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication practice table with a new set of problems.\n# Each filled row in the table (1, 5, 10, 15, 20) gets all five of its\n# cell values replaced with new \"NNN\u00d7N=\" expressions, by direct\n# (row, column) addressing so that values which collide with other\n# rows' old/new text cannot cross-contaminate each other.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Text = \"601\u00d75=\" },\n    @{ Row = 1;  Col = 2; Text = \"610\u00d73=\" },\n    @{ Row = 1;  Col = 3; Text = \"568\u00d77=\" },\n    @{ Row = 1;  Col = 4; Text = \"945\u00d78=\" },\n    @{ Row = 1;  Col = 5; Text = \"770\u00d72=\" },\n\n    @{ Row = 5;  Col = 1; Text = \"740\u00d78=\" },\n    @{ Row = 5;  Col = 2; Text = \"762\u00d78=\" },\n    @{ Row = 5;  Col = 3; Text = \"224\u00d78=\" },\n    @{ Row = 5;  Col = 4; Text = \"906\u00d73=\" },\n    @{ Row = 5;  Col = 5; Text = \"138\u00d76=\" },\n\n    @{ Row = 10; Col = 1; Text = \"762\u00d76=\" },\n    @{ Row = 10; Col = 2; Text = \"960\u00d77=\" },\n    @{ Row = 10; Col = 3; Text = \"434\u00d78=\" },\n    @{ Row = 10; Col = 4; Text = \"105\u00d74=\" },\n    @{ Row = 10; Col = 5; Text = \"513\u00d78=\" },\n\n    @{ Row = 15; Col = 1; Text = \"737\u00d78=\" },\n    @{ Row = 15; Col = 2; Text = \"346\u00d76=\" },\n    @{ Row = 15; Col = 3; Text = \"737\u00d73=\" },\n    @{ Row = 15; Col = 4; Text = \"742\u00d78=\" },\n    @{ Row = 15; Col = 5; Text = \"725\u00d77=\" },\n\n    @{ Row = 20; Col = 1; Text = \"816\u00d78=\" },\n    @{ Row = 20; Col = 2; Text = \"279\u00d78=\" },\n    @{ Row = 20; Col = 3; Text = \"701\u00d72=\" },\n    @{ Row = 20; Col = 4; Text = \"934\u00d76=\" },\n    @{ Row = 20; Col = 5; Text = \"652\u00d79=\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    $cellRange = $cell.Range\n    # Trim the trailing end-of-cell marker so only the visible text is replaced.\n    $cellRange.MoveEnd(12, -1) | Out-Null\n    $cellRange.Text = $u.Text\n}\n"}
